# Generate Report for Handoff
#
# The localization CI run finished: the "In Translation" status becomes
# "Ready for handoff" everywhere it appears (Overview rollup columns for
# zh-cn/de-de, plus the Status column on each per-language sheet), and the
# handoff/generate timestamps tick forward to reflect the new run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---------------------------------------------------
# Column E = zh-cn status, Column F = de-de status, Column G = Latest HO
# Xliff Generate Date.
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-31 08:48:01"

# --- zh-cn sheet --------------------------------------------------------
# Column C = Status, Column H = Latest Handoff Datetime.
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-31 08:47:55"

# --- de-de sheet --------------------------------------------------------
# Column C = Status, Column H = Latest Handoff Datetime.
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-31 08:48:01"

# --- Column widths --------------------------------------------------
# "Ready for handoff" is wider than "In Translation" was, so the Status
# columns widen to fit the new text.
$wsOverview.Columns.Item(5).ColumnWidth = 16.33
$wsOverview.Columns.Item(6).ColumnWidth = 16.33
$wsZhCn.Columns.Item(3).ColumnWidth = 16.33
$wsDeDe.Columns.Item(3).ColumnWidth = 16.33
